$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "60.774.36"
$ws.Range("E2").Value = "  -3.33%  "
$ws.Range("D3").Value = "2.913.44"
$ws.Range("E3").Value = "  -4.08%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("D6").Value = "144.04"
$ws.Range("E6").Value = "  -6.14%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "0.503"
$ws.Range("E8").Value = "  -2.96%  "
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "2.912.41"
$ws.Range("E9").Value = "  -3.92%  "
$ws.Range("D10").Value = "6.74"
$ws.Range("E10").Value = "  +6.15%  "
$ws.Range("E11").Value = "  -4.37%  "
$ws.Range("E12").Value = "  -3.77%  "
$ws.Range("E13").Value = "  -3.79%  "
$ws.Range("D14").Value = "33.45"
$ws.Range("E14").Value = "  -5.56%  "
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "3.396.50"
$ws.Range("E16").Value = "  -4.07%  "
$ws.Range("D17").Value = "60.725.08"
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("E18").Value = "  -4.64%  "
$ws.Range("D19").Value = "2.913.17"
$ws.Range("E19").Value = "  -4.05%  "
$ws.Range("D20").Value = "430.35"
$ws.Range("E20").Value = "  -4.69%  "
$ws.Range("D21").Value = "13.63"
$ws.Range("E21").Value = "  -4.42%  "
$ws.Range("E23").Value = "  -4.85%  "
$ws.Range("D24").Value = "80.41"
$ws.Range("E24").Value = "  -3.38%  "
$ws.Range("D25").Value = "10.90"
$ws.Range("E25").Value = "  -1.50%  "
$ws.Range("E26").Value = "  -4.77%  "
$ws.Range("D27").Value = "11.85"
$ws.Range("E27").Value = "  -3.23%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -4.67%  "
$ws.Range("E31").Value = "  -3.26%  "
$ws.Range("D32").Value = "2.15"
$ws.Range("E32").Value = "  -3.50%  "
$ws.Range("D33").Value = "26.51"
$ws.Range("E33").Value = "  -3.93%  "
$ws.Range("E34").Value = "  -4.61%  "
$ws.Range("D35").Value = "0.0₃0875"
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("E36").Value = "  -3.51%  "
$ws.Range("E37").Value = "  -4.83%  "
$ws.Range("D38").Value = "3.01"
$ws.Range("E38").Value = "  -4.91%  "
$ws.Range("D39").Value = "49.81"
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").Value = "2.00"
$ws.Range("E41").Value = "  -4.60%  "
$ws.Range("D42").Value = "8.69"
$ws.Range("E42").Value = "  -4.30%  "
$ws.Range("D43").Value = "0.295"
$ws.Range("E43").Value = "  -4.99%  "
$ws.Range("D44").Value = "41.68"
$ws.Range("E44").Value = "  -2.10%  "
$ws.Range("D45").Value = "377.53"
$ws.Range("E45").Value = "  -4.31%  "
$ws.Range("D46").Value = "0.0347"
$ws.Range("E46").Value = "  -3.56%  "
$ws.Range("D47").Value = "2.677.19"
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("D48").Value = "132.20"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D50").Value = "24.42"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("E51").Value = "  -2.03%  "
